$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 2021
$ws.Range("E8").Value = "hometown_lane.jpg"
$ws.Range("C8").Value = "Hometown Lane"
$ws.Range("B8").Value = "ふるさとこみち"
$ws.Range("D8").Value = "Shinkigensha"
$ws.Range("F8").Value = "supplement"

$ws.Range("A9").Select()
